$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 39, shifting existing rows 39-66 down to 40-67.
$ws.Rows.Item(39).Insert()

# Populate the newly inserted row 39 with the new weekly price record.
$ws.Range("A39").Value = 1
$ws.Range("B39").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C39").Value = "Arica y Parinacota"
$ws.Range("D39").Value = 44589
$ws.Range("E39").Value = 15
$ws.Range("F39").Value = 100112038
$ws.Range("G39").Value = "Cebollín baby"
$ws.Range("H39").Value = "Sin especificar"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 300
$ws.Range("K39").Value = 4500
$ws.Range("L39").Value = 5000
$ws.Range("M39").Value = 4750
$ws.Range("N39").Value = "$/paquete 1,5 a 2 kilos"
$ws.Range("O39").Value = "Región de Arica y Parinacota"
$ws.Range("P39").Value = 2375
$ws.Range("Q39").Value = 2
$ws.Range("R39").Value = "Hortaliza"
